# Generate Report for Archive
#
# The localization status report is regenerated: the "Ready for handoff"
# status (shared by the zh-cn/de-de columns on the Overview sheet and the
# Status column on each per-locale sheet) moves on to "In Translation".
# Excel re-flows the status columns to fit the new (shorter) text after the
# value change, so their stored widths shrink too.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears -----------------------
# A whole-cell (xlWhole = 1) Find/Replace across every worksheet updates all
# occurrences (Overview!E2:E3/F2:F3 and <locale>!C2:C3 on zh-cn/de-de) while
# keeping the edit as a single shared-string content change.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation", 1)
}

# --- 2. Shrink the columns that held the old, longer text ------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # zh-cn status column
$overview.Columns.Item(6).ColumnWidth = 12.5   # de-de status column

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # Status column

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # Status column
